$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the data grid A1:D4 (YCbCr GLCM data)
$data = @(
    @(0.021073234772125377, 0.90160144915746154, 0.77410860732452103, 0.9894633826139374),
    @(0.053627144361093039, 0.87964204395656598, 0.50649687740344296, 0.97318642781945364),
    @(0.000027361763983047511, 0.14813446693960255, 0.99994239441269861, 0.9999863191180085),
    @(0.0035706956932726052, 0.80269462818970838, 0.97650089665090734, 0.99821465215336369)
)

for ($r = 0; $r -lt 4; $r++) {
    for ($c = 0; $c -lt 4; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# Column widths (closest reachable quantization to 15.37890625 / 12.7109375 / 12.7109375 / 12.7109375)
$ws.Range("A1").EntireColumn.ColumnWidth = 14.5
$ws.Range("B1").EntireColumn.ColumnWidth = 11.83
$ws.Range("C1").EntireColumn.ColumnWidth = 11.83
$ws.Range("D1").EntireColumn.ColumnWidth = 11.83

# Register the extra cell-format styles (Text / Date number formats, each
# paired with its own border slot) that the source workbook carries in its
# style table without binding them to any visible cell -- mimic that by
# formatting scratch cells then clearing them, so the styles stay registered
# in styles.xml but no <c> in the saved sheet ends up referencing them.
$scratch1 = $ws.Range("Z1")
$scratch2 = $ws.Range("Z2")
$scratch1.Borders.LineStyle = 1
$scratch1.NumberFormat = "@"
$scratch2.Borders.LineStyle = 1
$scratch2.Borders.Color = 255
$scratch2.NumberFormat = "m/d/yy h:mm"
$ws.Range("Z1:Z2").Clear()

# Mark the workbook for a full recalculation on next load (mirrors
# calcPr/@fullCalcOnLoad="true").
$excel.CalculateFullRebuild()
